$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.130.43'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '3.257.56'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''397.60'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').Value = '''108.85'
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('D7').Value = '''0.578'
$ws.Range('E7').Value = '  +4.26%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.620'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').Value = '''39.26'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').Value = '''0.0955'
$ws.Range('E11').Value = '  +5.34%  '
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('D13').Value = '3.773.13'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').Value = '''8.26'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '''18.96'
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').Value = '3.264.43'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').Value = '''11.09'
$ws.Range('E18').Value = '  +3.90%  '
$ws.Range('D19').Value = '56.941.08'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = '''3.30'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').Value = '''0.0000108'
$ws.Range('E21').Value = '  +5.45%  '
$ws.Range('D22').Value = '''12.94'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '''294.55'
$ws.Range('E23').Value = '  -3.21%  '
$ws.Range('D24').Value = '''73.94'
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').Value = '''28.05'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('E27').Value = '  -3.85%  '
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').Value = '''7.43'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.113'
$ws.Range('E31').Value = '  +2.27%  '
$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').Value = '''11.16'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').Value = '''40.70'
$ws.Range('E34').Value = '  +12.34%  '
$ws.Range('D35').Value = '''0.0493'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('D37').Value = '''51.29'
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('D40').Value = '''3.00'
$ws.Range('E40').Value = '  -3.80%  '
$ws.Range('D41').Value = '''137.72'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '''0.284'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '''1.87'
$ws.Range('E44').Value = '  -2.88%  '
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('D46').Value = '''16.76'
$ws.Range('E46').Value = '  -2.72%  '
$ws.Range('D47').Value = '''22.36'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('E48').Value = '  +3.63%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '2.143.75'
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('E51').Value = '  -7.22%  '
